$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168 (pushes old rows 168..295 down to 169..296)
$ws.Rows("168:168").Insert()

# Populate the newly inserted row 168 with the new weekly record.
# (Same attributes as the existing "Puerro" / "Azul de Maquehue" / "Primera"
#  series for "Provincia de Cautín", but with a new date and new figures.)
$ws.Range("A168").Value = 10
$ws.Range("B168").Value = 'Vega Modelo de Temuco'
$ws.Range("C168").Value = 'La Araucanía'
$ws.Range("D168").Value = 45072
$ws.Range("E168").Value = 9
$ws.Range("F168").Value = 100112005
$ws.Range("G168").Value = 'Puerro'
$ws.Range("H168").Value = 'Azul de Maquehue'
$ws.Range("I168").Value = 'Primera'
$ws.Range("J168").Value = 40
$ws.Range("K168").Value = 10000
$ws.Range("L168").Value = 10000
$ws.Range("M168").Value = 10000
$ws.Range("N168").Value = '$/docena de paquetes'
$ws.Range("O168").Value = 'Provincia de Cautín'
$ws.Range("P168").Value = 833
$ws.Range("Q168").Value = 12
$ws.Range("R168").Value = 'Hortaliza'
